$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("N3").Value = 0
$ws.Range("H39").Value = 580.5
$ws.Range("I39").Value = 38.666668
$ws.Range("J39").Value = 1122.3334
$ws.Range("K39").Value = 116.000004
$ws.Range("L39").Value = 3367.0002
$ws.Range("M39").Value = 179.999996
$ws.Range("N39").Value = -3959.0002
$ws.Range("H42").Value = 44.375
$ws.Range("I42").Value = 38.75
$ws.Range("K42").Value = 116.25
$ws.Range("M42").Value = 113.75
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("M69").Value = 9000
$ws.Range("N69").Value = -10748
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3000
$ws.Range("K72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("M72").Value = 27000
$ws.Range("N72").Value = -35736
$ws.Range("H86").Value = 2116.5
$ws.Range("I86").Value = 2424.75
$ws.Range("K86").Value = 2424.75
$ws.Range("M86").Value = -1301.75
$ws.Range("H89").Value = 2116.5
$ws.Range("I89").Value = 2424.75
$ws.Range("K89").Value = 12123.75
$ws.Range("M89").Value = -6507.75
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H130").Value = 98497.25
$ws.Range("J130").Value = 98497.25
$ws.Range("L130").Value = 98497.25
$ws.Range("N130").Value = -108537.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2850
$ws.Range("I61").Value = 2850
$ws.Range("K61").Value = 2850
$ws.Range("M61").Value = -2638
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0
$ws.Range("H136").Value = 2850
$ws.Range("I136").Value = 2850
$ws.Range("K136").Value = 8550
$ws.Range("M136").Value = -6000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 503.33334
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 503.33334
$ws.Range("K11").Value = 0
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = 503.33334
$ws.Range("N11").Value = -783.33334
$ws.Range("H29").Value = 2062
$ws.Range("I29").Value = 103.666664
$ws.Range("J29").Value = 4999.5
$ws.Range("K29").Value = 103.666664
$ws.Range("L29").Value = 4999.5
$ws.Range("M29").Value = 185.333336
$ws.Range("N29").Value = -5577.5
$ws.Range("H86").Value = 1900
$ws.Range("I86").Value = 1200
$ws.Range("K86").Value = 1200
$ws.Range("M86").Value = -77
$ws.Range("H89").Value = 1900
$ws.Range("I89").Value = 1200
$ws.Range("K89").Value = 6000
$ws.Range("M89").Value = -384
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 7781.8335
$ws.Range("I134").Value = 6314
$ws.Range("K134").Value = 18942
$ws.Range("M134").Value = -16407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1500
$ws.Range("I58").Value = 1500
$ws.Range("K58").Value = 1500
$ws.Range("M58").Value = -1297
$ws.Range("H105").Value = 1249
$ws.Range("I105").Value = 1249
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1249
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = 498
$ws.Range("H122").Value = 1064.6666
$ws.Range("I122").Value = 997
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 2991
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -541
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 1599.4
$ws.Range("I132").Value = 1249.25
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3747.75
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1217.75
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 283.33334
$ws.Range("I86").Value = 225
$ws.Range("K86").Value = 675
$ws.Range("M86").Value = 511
$ws.Range("H89").Value = 283.33334
$ws.Range("I89").Value = 225
$ws.Range("K89").Value = 2025
$ws.Range("M89").Value = 3903
$ws.Range("H136").Value = 4753.5713
$ws.Range("I136").Value = 4753.5713
$ws.Range("K136").Value = 14260.7139
$ws.Range("M136").Value = -9160.713899999999
$ws.Range("H138").Value = 3797.2222
$ws.Range("I138").Value = 3709.375
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 11128.125
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = -5988.125
$ws.Range("N138").Value = -23780
$ws.Range("H139").Value = 4900
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 4900
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("M139").Value = 14700
$ws.Range("N139").Value = -24980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11100000
$ws.Range("I11").Value = 12500000
$ws.Range("K11").Value = 12500000
$ws.Range("M11").Value = -12499861
$ws.Range("H70").Value = 5038.1665
$ws.Range("I70").Value = 4945.8
$ws.Range("K70").Value = 4945.8
$ws.Range("M70").Value = -4675.8
$ws.Range("H73").Value = 5038.1665
$ws.Range("I73").Value = 4945.8
$ws.Range("K73").Value = 4945.8
$ws.Range("M73").Value = -4009.8
$ws.Range("H80").Value = 5500
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 5500
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H126").Value = 4978.4
$ws.Range("I126").Value = 4964
$ws.Range("K126").Value = 14892
$ws.Range("M126").Value = -12422
$ws.Range("H132").Value = 3774.889
$ws.Range("I132").Value = 3774.889
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11324.667
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8794.667000000001
$ws.Range("H134").Value = 99998.336
$ws.Range("J134").Value = 99998.336
$ws.Range("L134").Value = 299995.008
$ws.Range("N134").Value = -305065.008
$ws.Range("H136").Value = 31540.25
$ws.Range("J136").Value = 31540.25
$ws.Range("L136").Value = 94620.75
$ws.Range("N136").Value = -99720.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3112.8
$ws.Range("I16").Value = 3112.8
$ws.Range("K16").Value = 3112.8
$ws.Range("M16").Value = -2942.8
$ws.Range("H46").Value = 1670.2
$ws.Range("J46").Value = 1701
$ws.Range("L46").Value = 1701
$ws.Range("N46").Value = -2077
$ws.Range("H55").Value = 1569.7
$ws.Range("I55").Value = 1233.3334
$ws.Range("K55").Value = 1233.3334
$ws.Range("M55").Value = -1060.3334
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
$ws.Range("H128").Value = 89331.336
$ws.Range("J128").Value = 89331.336
$ws.Range("L128").Value = 89331.336
$ws.Range("N128").Value = -99291.336
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
$ws.Range("H136").Value = 1013.15
$ws.Range("I136").Value = 803.3158
$ws.Range("K136").Value = 2409.9474
$ws.Range("M136").Value = 140.0526
